$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Relative Share")

# Update row 8 (Year 2025) figures
$ws.Range("B8").Value = 0.7204610951008645
$ws.Range("C8").Value = 42.69932756964457
$ws.Range("E8").Value = 52.25744476464938
$ws.Range("F8").Value = 0.5763688760806917
$ws.Range("H8").Value = 1.633045148895293
